$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 3.4
$ws.Range("I3").Value = 2.2
$ws.Range("J3").Value = 4.33
$ws.Range("U3").Value = 3.75
$ws.Range("AA3").Value = 2.1
$ws.Range("AB3").Value = 1.67
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 41
$ws.Range("AG3").Value = 34
$ws.Range("AI3").Value = 7
$ws.Range("AK3").Value = 19
$ws.Range("AM3").Value = 6
$ws.Range("AN3").Value = 9.5
$ws.Range("AR3").Value = 41

# Row 5 updates
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 3.4
$ws.Range("L5").Value = 4
$ws.Range("AD5").Value = 9.5
$ws.Range("AF5").Value = 19
$ws.Range("AK5").Value = 15
$ws.Range("AM5").Value = 9.5
$ws.Range("AO5").Value = 12
